$d = $word.ActiveDocument

# Locate the "RG" attribute paragraph in the dictionary (the one whose
# explanatory text contains "carteira de identidade").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "carteira de identidade") {
        $target = $para
        break
    }
}

# 1) Indent this paragraph (adds <w:ind w:left="708"/> to its pPr).
$target.Format.LeftIndent = 35.4

# 2) Relocate the "_GoBack" bookmark so it sits right before the word that
#    follows "dos" inside this paragraph's explanatory text, splitting the
#    run in two.
$r = $target.Range.Duplicate
[void]$r.Find.Execute("funcion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $splitPoint)
